$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. RW03 rows (74-84): the "RegressionTest" column (E) flips from Yes to No,
#        since these are no longer the newest batch of test cases. ---
$ws.Range("E74:E84").Value = "No"

# --- 2. New rows 85-90 for the RW04 probate form feature. ---
# Seed formatting for the new rows from row 74 (closest fully-styled template row),
# then fix up column C (ScenarioName) to use the border+protection-only style
# (same as column B's style) instead of the wrap-text style used elsewhere.
$ws.Range("A74:E74").Copy()
$ws.Range("A85:E90").PasteSpecial(-4122)
$ws.Range("B74").Copy()
$ws.Range("C85:C90").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A85").Value = "TC_084"
$ws.Range("B85").Value = "probateFormsRW04.feature"
$ws.Range("C85").Value = "Open Estate"
$ws.Range("D85").Value = "Yes"
$ws.Range("E85").Value = "No"

$ws.Range("A86").Value = "TC_085"
$ws.Range("B86").Value = "probateFormsRW04.feature"
$ws.Range("C86").Value = "Verify, correct title is displayed on the form's header."
$ws.Range("D86").Value = "Yes"
$ws.Range("E86").Value = "No"

$ws.Range("A87").Value = "TC_086"
$ws.Range("B87").Value = "probateFormsRW04.feature"
$ws.Range("C87").Value = "Verify, county, and aka names are auto populated on the form."
$ws.Range("D87").Value = "Yes"
$ws.Range("E87").Value = "No"

$ws.Range("A88").Value = "TC_087"
$ws.Range("B88").Value = "probateFormsRW04.feature"
$ws.Range("C88").Value = "Verify, correct estate's name is displayed on the form."
$ws.Range("D88").Value = "Yes"
$ws.Range("E88").Value = "No"

$ws.Range("A89").Value = "TC_088"
$ws.Range("B89").Value = "probateFormsRW04.feature"
$ws.Range("C89").Value = "Verify, name of the decedent should be auto populated from the form."
$ws.Range("D89").Value = "Yes"
$ws.Range("E89").Value = "No"

$ws.Range("A90").Value = "TC_089"
$ws.Range("B90").Value = "probateFormsRW04.feature"
$ws.Range("C90").Value = "Verify, witnesses  name, address and signature should be editable and in yellow background."
$ws.Range("D90").Value = "Yes"
$ws.Range("E90").Value = "No"

# --- 3. Data validation (Yes/No list) for the newly added rows. ---
$ws.Range("D85:D90").Validation.Add(3, 1, 1, "Yes,No")
$ws.Range("E85:E90").Validation.Add(3, 1, 1, "Yes,No")

# --- 4. Update selection / active cell to match the new last data row. ---
$ws.Range("C85").Select()
